$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.219.58'
$ws.Range("E2").Value = '  -1.87%  '
$ws.Range("D3").Value = '1.582.24'
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.50'
$ws.Range("E6").Value = '  -2.40%  '
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("E8").Value = '  -1.48%  '
$ws.Range("E9").Value = '  -0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.50'
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '1.804.46'
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("D13").Value = '1.580.20'
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.43'
$ws.Range("E16").Value = '  -0.77%  '
$ws.Range("D17").Value = '26.212.88'
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '206.21'
$ws.Range("E21").Value = '  -1.94%  '
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("E23").Value = '  -2.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.86'
$ws.Range("E24").Value = '  -1.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.85'
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("E26").Value = '  -0.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.03'
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.20'
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0503'
$ws.Range("E30").Value = '  -1.47%  '
$ws.Range("E31").Value = '  -0.96%  '
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("D34").Value = '1.282.17'
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("E35").Value = '  +8.83%  '
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("E38").Value = '  -1.26%  '
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("E40").Value = '  -1.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.57'
$ws.Range("E41").Value = '  +3.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.768'
$ws.Range("E42").Value = '  -1.66%  '
$ws.Range("E43").Value = '  -2.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.18'
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("D45").Value = '1.717.64'
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.61'
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0508'
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₇0964'
$ws.Range("E50").Value = '  -9.32%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  -0.14%  '
